$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the naive QoQ error-series matrix (rows 2-24, cols B-K) with
# corrected values from the fixed naive component forecaster.

# Row 2
$ws.Cells.Item(2, 2).Value = 2.18989627706783
$ws.Cells.Item(2, 3).Value = 10.16349548557
$ws.Cells.Item(2, 4).Value = -7.724356307055293
$ws.Cells.Item(2, 5).Value = 0.4743409654747825
$ws.Cells.Item(2, 6).Value = 1.119950557408814
$ws.Cells.Item(2, 7).Value = -0.6780455854739627
$ws.Cells.Item(2, 8).Value = -1.635572794074753
$ws.Cells.Item(2, 9).Value = 1.198833335909301
$ws.Cells.Item(2, 10).Value = -0.4699347444161387
$ws.Cells.Item(2, 11).Value = 0.5504023208997705

# Row 3
$ws.Cells.Item(3, 2).Value = 9.247944187619927
$ws.Cells.Item(3, 3).Value = -8.279795221705264
$ws.Cells.Item(3, 4).Value = 0.1554212552770169
$ws.Cells.Item(3, 5).Value = 0.9107481729412195
$ws.Cells.Item(3, 6).Value = -0.8312630172827695
$ws.Cells.Item(3, 7).Value = -1.761027663680778
$ws.Cells.Item(3, 8).Value = 1.087260944115845
$ws.Cells.Item(3, 9).Value = -0.5745824983500689
$ws.Cells.Item(3, 10).Value = 0.4492110423354971
$ws.Cells.Item(3, 11).Value = -1.262772962148411

# Row 4
$ws.Cells.Item(4, 2).Value = -15.44502510429221
$ws.Cells.Item(4, 3).Value = -5.619233542269922
$ws.Cells.Item(4, 4).Value = -3.731309787990799
$ws.Cells.Item(4, 5).Value = -4.581790845914914
$ws.Cells.Item(4, 6).Value = -4.817112832495161
$ws.Cells.Item(4, 7).Value = -1.4324438270493
$ws.Cells.Item(4, 8).Value = -2.682060609296908
$ws.Cells.Item(4, 9).Value = -1.342494699581906
$ws.Cells.Item(4, 10).Value = -2.813101470050608
$ws.Cells.Item(4, 11).Value = 0.2004105337528701

# Row 5
$ws.Cells.Item(5, 2).Value = -3.939052395700697
$ws.Cells.Item(5, 3).Value = 3.684734974517905
$ws.Cells.Item(5, 4).Value = -2.333734529689821
$ws.Cells.Item(5, 5).Value = -0.835751817305987
$ws.Cells.Item(5, 6).Value = 0.5457380748606747
$ws.Cells.Item(5, 7).Value = -0.2700543632265874
$ws.Cells.Item(5, 8).Value = 0.2475156789190572
$ws.Cells.Item(5, 9).Value = -1.170528923740809
$ws.Cells.Item(5, 10).Value = 1.483733173460192
$ws.Cells.Item(5, 11).Value = -0.6491341121382845

# Row 6
$ws.Cells.Item(6, 2).Value = -0.04247155220008958
$ws.Cells.Item(6, 3).Value = -0.470863874926198
$ws.Cells.Item(6, 4).Value = -1.596660741476967
$ws.Cells.Item(6, 5).Value = 1.095969070482927
$ws.Cells.Item(6, 6).Value = -0.517307604149827
$ws.Cells.Item(6, 7).Value = 0.5080998195400784
$ws.Cells.Item(6, 8).Value = -1.214309943970497
$ws.Cells.Item(6, 9).Value = 1.615420670929393
$ws.Cells.Item(6, 10).Value = -0.6191041902945924
$ws.Cells.Item(6, 11).Value = 0.7318558435072805

# Row 7
$ws.Cells.Item(7, 2).Value = -0.062631421304428
$ws.Cells.Item(7, 3).Value = -1.345037697590573
$ws.Cells.Item(7, 4).Value = 1.014063016659791
$ws.Cells.Item(7, 5).Value = -0.4422155547685414
$ws.Cells.Item(7, 6).Value = 0.6140069800538022
$ws.Cells.Item(7, 7).Value = -1.155893297811842
$ws.Cells.Item(7, 8).Value = 1.681859354759485
$ws.Cells.Item(7, 9).Value = -0.5444222316837559
$ws.Cells.Item(7, 10).Value = 0.8021691354082593
$ws.Cells.Item(7, 11).Value = 0.1802518522544978

# Row 8
$ws.Cells.Item(8, 2).Value = -1.450536909439388
$ws.Cells.Item(8, 3).Value = 0.8879577326944336
$ws.Cells.Item(8, 4).Value = -0.2739674861353097
$ws.Cells.Item(8, 5).Value = 0.6746710726576846
$ws.Cells.Item(8, 6).Value = -1.144760458437209
$ws.Cells.Item(8, 7).Value = 1.744722106086416
$ws.Cells.Item(8, 8).Value = -0.4870707038900127
$ws.Cells.Item(8, 9).Value = 0.8466743124126384
$ws.Cells.Item(8, 10).Value = 0.2313257232809164
$ws.Cells.Item(8, 11).Value = 0.7341089884009673

# Row 9
$ws.Cells.Item(9, 2).Value = 0.2835171128466949
$ws.Cells.Item(9, 3).Value = -0.3381708634900031
$ws.Cells.Item(9, 4).Value = 0.9662183807149013
$ws.Cells.Item(9, 5).Value = -1.182467334523142
$ws.Cells.Item(9, 6).Value = 1.736138086012954
$ws.Cells.Item(9, 7).Value = -0.4033296262926443
$ws.Cells.Item(9, 8).Value = 0.8815014573406833
$ws.Cells.Item(9, 9).Value = 0.2578901063332134
$ws.Cells.Item(9, 10).Value = 0.7796078291475662
$ws.Cells.Item(9, 11).Value = 0.1037568572541728

# Row 10
$ws.Cells.Item(10, 2).Value = -0.4258116923289144
$ws.Cells.Item(10, 3).Value = 0.9137991786852327
$ws.Cells.Item(10, 4).Value = -1.123820937538137
$ws.Cells.Item(10, 5).Value = 1.741067231350056
$ws.Cells.Item(10, 6).Value = -0.4144229093415396
$ws.Cells.Item(10, 7).Value = 0.8942681565657697
$ws.Cells.Item(10, 8).Value = 0.2666196653805725
$ws.Cells.Item(10, 9).Value = 0.782109773296418
$ws.Cells.Item(10, 10).Value = 0.1100327947258539
$ws.Cells.Item(10, 11).Value = 0.5363482025840406

# Row 11
$ws.Cells.Item(11, 2).Value = 0.9362271583182413
$ws.Cells.Item(11, 3).Value = -1.121217397975688
$ws.Cells.Item(11, 4).Value = 1.718999545696041
$ws.Cells.Item(11, 5).Value = -0.4204048896458946
$ws.Cells.Item(11, 6).Value = 0.8889249063833586
$ws.Cells.Item(11, 7).Value = 0.2562511320102062
$ws.Cells.Item(11, 8).Value = 0.7736324511727942
$ws.Cells.Item(11, 9).Value = 0.1022604864190431
$ws.Cells.Item(11, 10).Value = 0.527724229051072
$ws.Cells.Item(11, 11).Value = 0.5878314504560218

# Row 12
$ws.Cells.Item(12, 2).Value = -1.099792826518468
$ws.Cells.Item(12, 3).Value = 1.839177394495253
$ws.Cells.Item(12, 4).Value = -0.5095992340596777
$ws.Cells.Item(12, 5).Value = 0.8591979025347809
$ws.Cells.Item(12, 6).Value = 0.2648772520401163
$ws.Cells.Item(12, 7).Value = 0.747692277747359
$ws.Cells.Item(12, 8).Value = 0.07948192339081561
$ws.Cells.Item(12, 9).Value = 0.5141061939024542
$ws.Cells.Item(12, 10).Value = 0.569310653270846
$ws.Cells.Item(12, 11).Value = 0.6926620895998143

# Row 13
$ws.Cells.Item(13, 2).Value = 1.796856762174863
$ws.Cells.Item(13, 3).Value = -0.5388037079154109
$ws.Cells.Item(13, 4).Value = 0.8561773027298597
$ws.Cells.Item(13, 5).Value = 0.2460159684286359
$ws.Cells.Item(13, 6).Value = 0.7273611466821548
$ws.Cells.Item(13, 7).Value = 0.06495176795656782
$ws.Cells.Item(13, 8).Value = 0.4974651484101862
$ws.Cells.Item(13, 9).Value = 0.5517268078420958
$ws.Cells.Item(13, 10).Value = 0.6761882060697715
$ws.Cells.Item(13, 11).Value = -0.2163793123768544

# Row 14
$ws.Cells.Item(14, 2).Value = -0.1986059672975008
$ws.Cells.Item(14, 3).Value = 0.9315483645137967
$ws.Cells.Item(14, 4).Value = 0.05492268211563681
$ws.Cells.Item(14, 5).Value = 0.7520550592688879
$ws.Cells.Item(14, 6).Value = 0.07433117151286378
$ws.Cells.Item(14, 7).Value = 0.4426543497987951
$ws.Cells.Item(14, 8).Value = 0.532832716953426
$ws.Cells.Item(14, 9).Value = 0.6613242052547258
$ws.Cells.Item(14, 10).Value = -0.2450910230062252
$ws.Cells.Item(14, 11).Value = 0.5297973106668776

# Row 15
$ws.Cells.Item(15, 2).Value = 1.382243236504047
$ws.Cells.Item(15, 3).Value = 0.1020333843279952
$ws.Cells.Item(15, 4).Value = 0.5117718284780797
$ws.Cells.Item(15, 5).Value = 0.1102762525787266
$ws.Cells.Item(15, 6).Value = 0.4406035923592023
$ws.Cells.Item(15, 7).Value = 0.4562644775173959
$ws.Cells.Item(15, 8).Value = 0.6339248108794424
$ws.Cells.Item(15, 9).Value = -0.2713800732039505
$ws.Cells.Item(15, 10).Value = 0.4864000199535451
$ws.Cells.Item(15, 11).Value = 0.2305062539156956

# Row 16
$ws.Cells.Item(16, 2).Value = 0.4143589712515336
$ws.Cells.Item(16, 3).Value = 0.645060568213604
$ws.Cells.Item(16, 4).Value = -0.06880216923746499
$ws.Cells.Item(16, 5).Value = 0.4687978311297695
$ws.Cells.Item(16, 6).Value = 0.4922897887097893
$ws.Cells.Item(16, 7).Value = 0.5963832406752287
$ws.Cells.Item(16, 8).Value = -0.2769224122913118
$ws.Cells.Item(16, 9).Value = 0.4905753488922937
$ws.Cells.Item(16, 10).Value = 0.2193215401759246

# Row 17
$ws.Cells.Item(17, 2).Value = 0.8806510593214452
$ws.Cells.Item(17, 3).Value = 0.01587657163684458
$ws.Cells.Item(17, 4).Value = 0.3231974644960524
$ws.Cells.Item(17, 5).Value = 0.5022435315690319
$ws.Cells.Item(17, 6).Value = 0.6117484103747737
$ws.Cells.Item(17, 7).Value = -0.3157707472487348
$ws.Cells.Item(17, 8).Value = 0.4747295398651078
$ws.Cells.Item(17, 9).Value = 0.2109873117084238

# Row 18
$ws.Cells.Item(18, 2).Value = 0.3269007581182274
$ws.Cells.Item(18, 3).Value = 0.4402862389072476
$ws.Cells.Item(18, 4).Value = 0.3397793901194651
$ws.Cells.Item(18, 5).Value = 0.6407073608399754
$ws.Cells.Item(18, 6).Value = -0.2800278747938312
$ws.Cells.Item(18, 7).Value = 0.4434018934814807
$ws.Cells.Item(18, 8).Value = 0.2077622620068982

# Row 19
$ws.Cells.Item(19, 2).Value = 0.687111416311362
$ws.Cells.Item(19, 3).Value = 0.3574352576573546
$ws.Cells.Item(19, 4).Value = 0.5460851217289808
$ws.Cells.Item(19, 5).Value = -0.2478104865453511
$ws.Cells.Item(19, 6).Value = 0.455790019309198
$ws.Cells.Item(19, 7).Value = 0.1868984584576193

# Row 20
$ws.Cells.Item(20, 2).Value = 0.5971116170735665
$ws.Cells.Item(20, 3).Value = 0.6311093651981943
$ws.Cells.Item(20, 4).Value = -0.3642950649033653
$ws.Cells.Item(20, 5).Value = 0.4708952909610505
$ws.Cells.Item(20, 6).Value = 0.2101374940836094

# Row 21
$ws.Cells.Item(21, 2).Value = 0.7959694422322816
$ws.Cells.Item(21, 3).Value = -0.350757792655449
$ws.Cells.Item(21, 4).Value = 0.4079223719710875
$ws.Cells.Item(21, 5).Value = 0.2201756597651073

# Row 22
$ws.Cells.Item(22, 2).Value = -0.09690875079004102
$ws.Cells.Item(22, 3).Value = 0.5072404221531239
$ws.Cells.Item(22, 4).Value = 0.1085991175498651

# Row 23
$ws.Cells.Item(23, 2).Value = 0.5513001133925729
$ws.Cells.Item(23, 3).Value = 0.130019622424466

# Row 24
$ws.Cells.Item(24, 2).Value = 0.3662627537369125

